$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing column E ("First Name"), making room
# for the new "Send Confirmation Email" column right after "Email".
$ws.Columns.Item(5).Insert()

# The new column sits next to the (wide) Email column, so give it the same
# display width instead of leaving it at the sheet default.
$ws.Columns.Item(5).ColumnWidth = 49.836666666666666

# Header for the newly inserted column.
$ws.Range("E1").Value = "Send Confirmation Email"

# Populate the new column: "No" for the Series A rows, "Yes" for the Pool 1 rows.
$ws.Range("E2").Value = "No"
$ws.Range("E3").Value = "No"
$ws.Range("E4").Value = "No"
$ws.Range("E5").Value = "No"
$ws.Range("E6").Value = "Yes"
$ws.Range("E7").Value = "Yes"

# Update the email domain across the sheet.
$ws.Cells.Replace("mycompany.com", "myfirm.com")

# Drop the mailto hyperlinks that used to live on the Email column.
$ws.Hyperlinks.Delete()

# Match the author's final selection.
$ws.Range("D2:D7").Select()
